$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.750.24'
$ws.Range("E2").Value = '  -0.88%  '

Set-TextValue $ws.Range("D3") '2.349.00'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("E4").Value = '  +0.09%  '

Set-TextValue $ws.Range("D5") '0.675'
$ws.Range("E5").Value = '  -0.24%  '

Set-TextValue $ws.Range("D6") '238.85'
$ws.Range("E6").Value = '  -0.31%  '

Set-TextValue $ws.Range("D7") '73.56'
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  -0.03%  '

Set-TextValue $ws.Range("D9") '0.605'
$ws.Range("E9").Value = '  +8.42%  '

Set-TextValue $ws.Range("D10") '0.100'
$ws.Range("E10").Value = '  -1.59%  '

Set-TextValue $ws.Range("D11") '58.37'
$ws.Range("E11").Value = '  +1.23%  '

Set-TextValue $ws.Range("D12") '33.24'
$ws.Range("E12").Value = '  +5.86%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D13") '7.33'
$ws.Range("E13").Value = '  +0.18%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D14") '0.108'
$ws.Range("E14").Value = '  +0.21%  '

Set-TextValue $ws.Range("D15") '2.698.08'
$ws.Range("E15").Value = '  -0.32%  '

Set-TextValue $ws.Range("D16") '16.40'
$ws.Range("E16").Value = '  -2.49%  '

Set-TextValue $ws.Range("D17") '0.902'
$ws.Range("E17").Value = '  -0.95%  '

Set-TextValue $ws.Range("D18") '2.357.70'
$ws.Range("E18").Value = '  -0.38%  '

Set-TextValue $ws.Range("D19") '43.657.62'
$ws.Range("E19").Value = '  -1.07%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D20") '0.0000101'
$ws.Range("E20").Value = '  -0.55%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D21") '6.77'
$ws.Range("E21").Value = '  +2.64%  '

Set-TextValue $ws.Range("D22") '77.12'
$ws.Range("E22").Value = '  -1.27%  '

Set-TextValue $ws.Range("D23") '256.80'
$ws.Range("E23").Value = '  +0.70%  '

Set-TextValue $ws.Range("D24") '1.97'
$ws.Range("E24").Value = '  +24.03%  '

Set-TextValue $ws.Range("D25") '0.999'
$ws.Range("E25").Value = '  -0.06%  '

Set-TextValue $ws.Range("D26") '3.72'
$ws.Range("E26").Value = '  -1.06%  '

Set-TextValue $ws.Range("D27") '2.48'
$ws.Range("E27").Value = '  -1.51%  '

Set-TextValue $ws.Range("D28") '10.58'
$ws.Range("E28").Value = '  -0.51%  '

Set-TextValue $ws.Range("D29") '2.26'
$ws.Range("E29").Value = '  -1.76%  '

Set-TextValue $ws.Range("D30") '22.77'
$ws.Range("E30").Value = '  +1.22%  '

Set-TextValue $ws.Range("D31") '175.73'
$ws.Range("E31").Value = '  +1.03%  '

Set-TextValue $ws.Range("D32") '0.129'
$ws.Range("E32").Value = '  -2.25%  '

Set-TextValue $ws.Range("D33") '0.136'
$ws.Range("E33").Value = '  +2.90%  '

Set-TextValue $ws.Range("D34") '0.0761'
$ws.Range("E34").Value = '  +2.30%  '

Set-TextValue $ws.Range("D35") '5.49'
$ws.Range("E35").Value = '  +4.36%  '

Set-TextValue $ws.Range("D36") '5.17'
$ws.Range("E36").Value = '  -2.34%  '

Set-TextValue $ws.Range("D37") '3.77'
$ws.Range("E37").Value = '  -2.40%  '

Set-TextValue $ws.Range("D38") '2.35'
$ws.Range("E38").Value = '  -3.69%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D39") '0.0281'
$ws.Range("E39").Value = '  +3.08%  '

$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D40") '6.22'
$ws.Range("E40").Value = '  -4.55%  '

Set-TextValue $ws.Range("D41") '0.111'
$ws.Range("E41").Value = '  +11.73%  '

Set-TextValue $ws.Range("D42") '0.206'
$ws.Range("E42").Value = '  +10.21%  '

$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range("D43") '61.86'
$ws.Range("E43").Value = '  +17.55%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D44") '8.99'
$ws.Range("E44").Value = '  +0.61%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D45") '18.84'
$ws.Range("E45").Value = '  -1.92%  '

$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range("D46") '4.67'
$ws.Range("E46").Value = '  +4.24%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D47") '2.49'
$ws.Range("E47").Value = '  +1.79%  '

$ws.Range("B48").Value = 'BinanceUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D48") '1.00'
$ws.Range("E48").Value = '  +0.06%  '

$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D49") '1.23'
$ws.Range("E49").Value = '  -1.76%  '

Set-TextValue $ws.Range("D50") '99.69'
$ws.Range("E50").Value = '  -0.37%  '

Set-TextValue $ws.Range("D51") '1.15'
$ws.Range("E51").Value = '  -1.11%  '
